$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New rows 30-33 appended to the competitor list
#   Columns: A=Competition B=Category C=Website D=State
#            H=Platforms I=Business Model J=Notes
# ---------------------------------------------------------------------------

# Row 30 - Omnifocus
$ws.Range("A30").Value = "Omnifocus"
$ws.Range("B30").Value = "Task management"
$ws.Range("C30").Value = "http://www.omnigroup.com/products/omnifocus/"
$ws.Range("D30").Value = "GA"
$ws.Range("H30").Value = "Mac, iPhone, iPad"
$ws.Range("I30").Value = "`$80 license"
$ws.Range("J30").Value = "cloud sync, e-mail integration"

# Row 31 - Daytum (has a hyperlink on the Website cell)
$ws.Range("A31").Value = "Daytum"
$ws.Range("B31").Value = "lists & statistics"
$ws.Range("C31").Value = "http://www.daytum.com/"
$ws.Range("D31").Value = "GA"
$ws.Hyperlinks.Add($ws.Range("C31"), "http://www.daytum.com/")
$ws.Range("C31").Style = "Hyperlink"

# Row 32 - ReQall
$ws.Range("A32").Value = "ReQall"
$ws.Range("B32").Value = "Voice-enabled memory aid"
$ws.Range("C32").Value = "http://www.reqall.com/about"
$ws.Range("D32").Value = "GA"
$ws.Range("H32").Value = "iPhone, Android, BB"
$ws.Range("I32").Value = "Freemium (25/yr)"
$ws.Range("J32").Value = "connectors for outlook, google cal, evernote"

# Row 33 - Carbonfin outliner
$ws.Range("A33").Value = "Carbonfin outliner"
$ws.Range("B33").Value = "Task management"
$ws.Range("C33").Value = "http://carbonfin.com/"
$ws.Range("D33").Value = "GA"
$ws.Range("H33").Value = "iPhone, iPad"
$ws.Range("I33").Value = "license fee on device"
$ws.Range("J33").Value = "sync, collab between users"

# ---------------------------------------------------------------------------
# Extend the table / autofilter to cover the new rows
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J33"))

# ---------------------------------------------------------------------------
# Column I (Business Model) widened to fit the new content
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 16.67

# ---------------------------------------------------------------------------
# Selection / scroll position moves to the newly added area
# ---------------------------------------------------------------------------
$null = $ws.Range("A32").Select()
